# "Pedoe picked 3 questions"
# - Marks the three existing Linked-List rows whose PIC is "Pedoe" (rows 2,4,5
#   on the "Easy" sheet) with Language = Javascript.
# - Adds three new Stack & Queue rows (Valid Parentheses / Next Greater
#   Element II / Decoding String) worked on by Pedoe, Status = On-going,
#   Language = Javascript.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Easy")

# --- 1. Copy the "Language" column formatting onto the three existing cells
#        that don't have it yet (F2, F4, F5 currently carry the blank/default
#        style used for rows with no language picked).
$ws.Range("F3").Copy()
$ws.Range("F2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("F5").PasteSpecial(-4122)

# --- 2. Prime the formatting of the new rows (14-16) off row 11, which is
#        the first row of the "Stack & Queue" block.
$ws.Range("E11").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("F11").Copy()
$ws.Range("F14").PasteSpecial(-4122)

$ws.Range("E11").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F11").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("E11").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F11").Copy()
$ws.Range("F16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 3. Fill in the "Javascript" language value everywhere it is newly used
#        (existing rows first, so it lands in the shared-string table before
#        the new question titles below).
$ws.Range("F2").Value = "Javascript"
$ws.Range("F4").Value = "Javascript"
$ws.Range("F5").Value = "Javascript"
$ws.Range("F14").Value = "Javascript"
$ws.Range("F15").Value = "Javascript"
$ws.Range("F16").Value = "Javascript"

# --- 4. New question rows.
$ws.Range("A14").Value = 20
$ws.Range("B14").Value = "Valid Parentheses"
$ws.Range("C14").Value = "Stack & Queue"
$ws.Range("D14").Value = "Pedoe"
$ws.Range("E14").Value = "On-going"

$ws.Range("A15").Value = 503
$ws.Range("B15").Value = "Next Greater Element II"
$ws.Range("C15").Value = "Stack & Queue"
$ws.Range("D15").Value = "Pedoe"
$ws.Range("E15").Value = "On-going"

$ws.Range("A16").Value = 394
$ws.Range("B16").Value = "Decoding String"
$ws.Range("C16").Value = "Stack & Queue"
$ws.Range("D16").Value = "Pedoe"
$ws.Range("E16").Value = "On-going"

# --- 5. Leave the selection where Excel would after typing the last row.
[void]$ws.Range("E16").Select()
